$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '61.839.49'
Set-TextCell 'E2' '  +0.79%  '
Set-TextCell 'D3' '3.424.04'
Set-TextCell 'E3' '  +1.07%  '
Set-TextCell 'D4' '1.00'
Set-TextCell 'D5' '408.99'
Set-TextCell 'E5' '  +1.45%  '
Set-TextCell 'D6' '128.36'
Set-TextCell 'E6' '  -1.11%  '
Set-TextCell 'E7' '  +7.59%  '
Set-TextCell 'D8' '1.00'
Set-TextCell 'E8' '  -0.03%  '
Set-TextCell 'D9' '0.734'
Set-TextCell 'E9' '  +8.34%  '
Set-TextCell 'D10' '0.139'
Set-TextCell 'E10' '  +10.33%  '
Set-TextCell 'D11' '42.58'
Set-TextCell 'E11' '  +2.66%  '
Set-TextCell 'D12' '9.13'
Set-TextCell 'E12' '  +10.04%  '
Set-TextCell 'E13' '  +0.21%  '
Set-TextCell 'D14' '3.964.14'
Set-TextCell 'E14' '  +0.70%  '
Set-TextCell 'D15' '21.25'
Set-TextCell 'E15' '  +7.90%  '
Set-TextCell 'D16' '0.0000205'
Set-TextCell 'E16' '  +45.56%  '
Set-TextCell 'D17' '3.412.44'
Set-TextCell 'E17' '  -0.02%  '
Set-TextCell 'D18' '12.40'
Set-TextCell 'E18' '  +5.83%  '
Set-TextCell 'E19' '  +7.48%  '
Set-TextCell 'D20' '61.907.07'
Set-TextCell 'E20' '  +0.90%  '
Set-TextCell 'D21' '444.68'
Set-TextCell 'E21' '  +43.47%  '
Set-TextCell 'D22' '91.37'
Set-TextCell 'E22' '  +10.42%  '
Set-TextCell 'D23' '3.19'
Set-TextCell 'E23' '  +1.59%  '
Set-TextCell 'D24' '12.94'
Set-TextCell 'E24' '  +2.20%  '
Set-TextCell 'D25' '3.24'
Set-TextCell 'E25' '  +3.70%  '
Set-TextCell 'D26' '33.05'
Set-TextCell 'E26' '  +12.72%  '
Set-TextCell 'D27' '8.69'
Set-TextCell 'E27' '  +8.64%  '
Set-TextCell 'E28' '  -0.37%  '
Set-TextCell 'E29' '  +2.48%  '
Set-TextCell 'D30' '7.61'
Set-TextCell 'E30' '  -6.87%  '
Set-TextCell 'D31' '11.96'
Set-TextCell 'E31' '  +6.51%  '
Set-TextCell 'E32' '  +0.01%  '
Set-TextCell 'E33' '  +0.30%  '
Set-TextCell 'D34' '42.66'
Set-TextCell 'E34' '  -1.88%  '
Set-TextCell 'E35' '  -0.04%  '
Set-TextCell 'D36' '0.0499'
Set-TextCell 'E36' '  +4.10%  '
Set-TextCell 'D37' '53.20'
Set-TextCell 'E37' '  +3.84%  '
Set-TextCell 'D38' '0.999'
Set-TextCell 'E38' '  -0.18%  '
Set-TextCell 'D39' '3.39'
Set-TextCell 'E39' '  +2.01%  '
Set-TextCell 'E40' '  +8.21%  '
Set-TextCell 'D41' '2.93'
Set-TextCell 'E41' '  +0.42%  '
Set-TextCell 'B42' 'Monero'
Set-TextCell 'C42' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D42' '142.55'
Set-TextCell 'E42' '  +2.90%  '
Set-TextCell 'B43' 'TheGraph'
Set-TextCell 'C43' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell 'D43' '0.310'
Set-TextCell 'E43' '  -3.40%  '
Set-TextCell 'D44' '4.24'
Set-TextCell 'E44' '  +8.68%  '
Set-TextCell 'D45' '1.99'
Set-TextCell 'E45' '  +1.87%  '
Set-TextCell 'D46' '2.51'
Set-TextCell 'E46' '  +13.52%  '
Set-TextCell 'D47' '16.56'
Set-TextCell 'E47' '  +0.38%  '
Set-TextCell 'D48' '22.42'
Set-TextCell 'E48' '  +6.83%  '
Set-TextCell 'B49' 'ThetaToken'
Set-TextCell 'C49' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell 'D49' '2.14'
Set-TextCell 'E49' '  +12.83%  '
Set-TextCell 'B50' 'RocketPoolETH'
Set-TextCell 'C50' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 'D50' '3.772.64'
Set-TextCell 'E50' '  +1.04%  '
Set-TextCell 'D51' '2.129.83'
Set-TextCell 'E51' '  +2.02%  '
